$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D11").Value = 17.72164948453608
$ws.Range("E11").Value = 4147010.30927835
$ws.Range("F11").Value = 241647.2934758712
$ws.Range("H11").Value = 97
